$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.04712155826069875
    "D2" = 0.03749071308083529
    "E2" = 0.1171343112676517
    "F2" = 3.62692846977319
    "G2" = 0.00258229471374426
    "I2" = 0.7307442631768311
    "J2" = 0.2319293090669419
    "K2" = 2.87563779270198
    "N2" = 2.247421581053047
    "B3" = 0.04112866080940591
    "D3" = 0.03671984399024097
    "E3" = 0.1148255713462412
    "F3" = 3.596872982392142
    "G3" = 0.002588201780499905
    "I3" = 0.7347034126563621
    "J3" = 0.2273828962592361
    "K3" = 2.737902908737397
    "N3" = 2.266419069499172
    "B4" = 0.03744387610867506
    "D4" = 0.03626107055945837
    "E4" = 0.1134734875882231
    "F4" = 3.580478423878503
    "G4" = 0.002592017461929164
    "I4" = 0.7374236260298979
    "J4" = 0.2247328354698439
    "K4" = 2.655144178472028
    "N4" = 2.278794512233077
    "B5" = 0.03594115326985303
    "D5" = 0.03607779353138696
    "E5" = 0.1129389259826468
    "F5" = 3.574313417647332
    "G5" = 0.002593620008497138
    "I5" = 0.7386047821500021
    "J5" = 0.2236883131504968
    "K5" = 2.621871599257531
    "N5" = 2.284016057074517
    "B6" = 0.03569156251328565
    "D6" = 0.0360475831874254
    "E6" = 0.1128511530896432
    "F6" = 3.573320824836756
    "G6" = 0.002593888991503633
    "I6" = 0.7388052971835535
    "J6" = 0.223517004085835
    "K6" = 2.616373938997697
    "N6" = 2.284893859601745
    "B7" = 0.0374236142988309
    "D7" = 0.03625858390687853
    "E7" = 0.113466211862729
    "F7" = 3.58039319422484
    "G7" = 0.002592038881304178
    "I7" = 0.7374392615007608
    "J7" = 0.2247186055699402
    "K7" = 2.654693625703715
    "N7" = 2.27886420970755
    "B8" = 0.04505635486707149
    "D8" = 0.03722191041453726
    "E8" = 0.1163246364571577
    "F8" = 3.616136628722529
    "G8" = 0.00258429240471747
    "I8" = 0.7320492803848033
    "J8" = 0.2303322412079751
    "K8" = 2.827768975093136
    "N8" = 2.253824178704683
    "B9" = 0.0599778826433095
    "D9" = 0.03922566229568503
    "E9" = 0.1224521667354637
    "F9" = 3.702671252880009
    "G9" = 0.002570591183748227
    "I9" = 0.7237793448119945
    "J9" = 0.2424716202640695
    "K9" = 3.181702885845425
    "N9" = 2.210375083417212
    "B10" = 0.07090606800396415
    "D10" = 0.04076697340690316
    "E10" = 0.1272765922584114
    "F10" = 3.776423052490031
    "G10" = 0.002561422050634579
    "I10" = 0.7191120078905087
    "J10" = 0.2520933448739981
    "K10" = 3.45085812704076
    "N10" = 2.181918199307873
    "B11" = 0.07586876236537421
    "D11" = 0.0414830335590608
    "E11" = 0.1295423566324203
    "F11" = 3.812218152760181
    "G11" = 0.00255744324178289
    "I11" = 0.7172960478903434
    "J11" = 0.256626178740234
    "K11" = 3.575343253676579
    "N11" = 2.169728869911239
    "B12" = 0.07774663937131265
    "D12" = 0.04175631590188544
    "E12" = 0.130410636432309
    "F12" = 3.82609827129653
    "G12" = 0.002555964039822513
    "I12" = 0.7166526927366021
    "J12" = 0.2583652842601936
    "K12" = 3.622780883032817
    "N12" = 2.165222154901699
    "B13" = 0.07734226908647202
    "D13" = 0.04169736539344626
    "E13" = 0.1302231785057799
    "F13" = 3.823094432352462
    "G13" = 0.002556281392346597
    "I13" = 0.7167892780181546
    "J13" = 0.2579897274142127
    "K13" = 3.612551035021625
    "N13" = 2.166187896560878
    "B14" = 0.07602328511536882
    "D14" = 0.04150547415792261
    "E14" = 0.1296135840850638
    "F14" = 3.813353544902554
    "G14" = 0.00255732099692807
    "I14" = 0.7172422298980621
    "J14" = 0.256768801465455
    "K14" = 3.579239985657182
    "N14" = 2.169355910350433
    "B15" = 0.07521518413676631
    "D15" = 0.04138821156119832
    "E15" = 0.1292415309525481
    "F15" = 3.807429406323934
    "G15" = 0.002557961359472719
    "I15" = 0.7175254502829702
    "J15" = 0.2560239010329326
    "K15" = 3.558874902568107
    "N15" = 2.171310633732141
    "B16" = 0.07058155833856006
    "D16" = 0.04072047565712467
    "E16" = 0.1271299550928902
    "F16" = 3.774129143326149
    "G16" = 0.00256168593043007
    "I16" = 0.7192368785207677
    "J16" = 0.2518002674646453
    "K16" = 3.442764200529325
    "N16" = 2.182730046771724
    "B17" = 0.06773666984432225
    "D17" = 0.04031464706141463
    "E17" = 0.1258528255597327
    "F17" = 3.754277251288357
    "G17" = 0.002564019964700516
    "I17" = 0.7203655660244728
    "J17" = 0.2492492711419345
    "K17" = 3.372060465889945
    "N17" = 2.189929380688667
    "B18" = 0.06609956554437701
    "D18" = 0.04008262964060805
    "E18" = 0.1251249439334785
    "F18" = 3.743070131670748
    "G18" = 0.002565380546935891
    "I18" = 0.7210436677777672
    "J18" = 0.2477966634786384
    "K18" = 3.331585809656076
    "N18" = 2.194141365143224
    "B19" = 0.06554513746246471
    "D19" = 0.04000431430452522
    "E19" = 0.1248796428934185
    "F19" = 3.73931178697859
    "G19" = 0.002565844331899619
    "I19" = 0.7212782228786629
    "J19" = 0.2473073467683093
    "K19" = 3.317914695182651
    "N19" = 2.195579672339861
    "B20" = 0.06803959695052697
    "D20" = 0.04035770298534658
    "E20" = 0.1259880854218096
    "F20" = 3.756368646814849
    "G20" = 0.002563769629721158
    "I20" = 0.7202424222501094
    "J20" = 0.2495193105291662
    "K20" = 3.379567073973533
    "N20" = 2.18915563578193
    "B21" = 0.07641074154835792
    "D21" = 0.04156177971091779
    "E21" = 0.1297923571621666
    "F21" = 3.816205831535257
    "G21" = 0.002557014895167452
    "I21" = 0.7171079833618847
    "J21" = 0.2571268014508092
    "K21" = 3.589016135545535
    "N21" = 2.168422422969428
    "B22" = 0.08187361026136841
    "D22" = 0.04236109677405153
    "E22" = 0.1323386402147477
    "F22" = 3.857210110998551
    "G22" = 0.002552760423443166
    "I22" = 0.7153177775114514
    "J22" = 0.2622306805655938
    "K22" = 3.727640886175493
    "N22" = 2.155508307133303
    "B23" = 0.07895877092067849
    "D23" = 0.04193335900810524
    "E23" = 0.1309741337873547
    "F23" = 3.835150934702199
    "G23" = 0.002555016515974863
    "I23" = 0.7162495634215063
    "J23" = 0.2594945018947925
    "K23" = 3.653493992843664
    "N23" = 2.162342453936958
    "B24" = 0.06790264839078475
    "D24" = 0.04033823335471709
    "E24" = 0.1259269146415107
    "F24" = 3.755422485446047
    "G24" = 0.002563882747859619
    "I24" = 0.720298004595449
    "J24" = 0.2493971821493659
    "K24" = 3.376172795196624
    "N24" = 2.189505218534372
    "B25" = 0.05594683759161967
    "D25" = 0.03867139743773151
    "E25" = 0.1207381775502157
    "F25" = 3.677485644776255
    "G25" = 0.002574139388915004
    "I25" = 0.7257697004243369
    "J25" = 0.2390651383072964
    "K25" = 3.08437187981508
    "N25" = 2.221522036067384
}

foreach ($key in $values.Keys) {
    $ws.Range($key).Value = $values[$key]
}
